$d = $word.ActiveDocument

# Locate the paragraph that ends the "Requisitos" section
# ("LOQ4233: Gestão de Negócios (Requisito fraco)") so we can find the
# block of footer-like paragraphs that follow it (an empty paragraph,
# "Ver no Jupiter Salvar em pdf Salvar em docx", and the copyright
# notice) and remove them.
$anchorStart = $d.Content.Duplicate
$foundStart = $anchorStart.Find.Execute(
    "LOQ4233: Gestão de Negócios (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorEnd = $d.Content.Duplicate
$foundEnd = $anchorEnd.Find.Execute(
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundStart -and $foundEnd) {
    # Range starting right after the "LOQ4233..." paragraph mark, through
    # the end of the copyright paragraph (including its paragraph mark),
    # so the whole block of 3 paragraphs disappears.
    $deleteStart = $anchorStart.Paragraphs(1).Next().Range.Start
    $deleteEnd = $anchorEnd.Paragraphs(1).Range.End

    $toDelete = $d.Range($deleteStart, $deleteEnd)
    $toDelete.Delete()
}
